$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("21:21").Delete()
